# Insert a new "Emission" column into the Externalities sheet (between Tech
# and External Cost), shifting the existing External Cost / Mode_Of_Operation /
# EmissionActivityRatio / EmissionsPenalty / Final Unit columns one to the
# right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Externalities")

$ws.Columns("B").Insert()
$ws.Range("B1").Value = "Emission"
$ws.Columns("B").ColumnWidth = 11.17

$ws.Activate()
$ws.Range("B9").Select()
